$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 268, shifting existing rows 268-313 down to 269-314
$ws.Rows.Item(268).Insert()

# Fill the new row 268 with the new data
$ws.Cells.Item(268, 1).Value = 8
$ws.Cells.Item(268, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(268, 3).Value = "Coquimbo"
$ws.Cells.Item(268, 4).Value = 45218
$ws.Cells.Item(268, 5).Value = 4
$ws.Cells.Item(268, 6).Value = 100112001
$ws.Cells.Item(268, 7).Value = "Berenjena"
$ws.Cells.Item(268, 8).Value = "Sin especificar"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 400
$ws.Cells.Item(268, 11).Value = 7500
$ws.Cells.Item(268, 12).Value = 8000
$ws.Cells.Item(268, 13).Value = 7750
$ws.Cells.Item(268, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(268, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(268, 16).Value = 155
$ws.Cells.Item(268, 17).Value = 50
$ws.Cells.Item(268, 18).Value = "Hortaliza"
